$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the daily price log for
# "Ciboulette" at Mercado Mayorista Lo Valledor de Santiago. Insert a
# whole new row before the existing row 294 (shifting every row below it
# down by one, extending the used range from A1:R390 to A1:R391), then
# populate the new row with the new record's data.
$ws.Rows("294").Insert()

$ws.Range("A294").Value = 6
$ws.Range("B294").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C294").Value = "Metropolitana"
$ws.Range("D294").Value = 44627
$ws.Range("E294").Value = 13
$ws.Range("F294").Value = 100112039
$ws.Range("G294").Value = "Ciboulette"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 570
$ws.Range("K294").Value = 1400
$ws.Range("L294").Value = 1500
$ws.Range("M294").Value = 1447
$ws.Range("N294").Value = "$/docena de atados"
$ws.Range("O294").Value = "Región Metropolitana"
$ws.Range("P294").Value = 482
$ws.Range("Q294").Value = 3
$ws.Range("R294").Value = "Hortaliza"
